# Insert a new weekly price record before row 91 (Fecha 2021-09-29, the
# most recent week), pushing the existing history for rows 91-128 down by
# one row (92-129) and bumping the sheet's used range accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("91:91").Insert()

$ws.Range("A91").Value = 3
$ws.Range("B91").Value = "Femacal de La Calera"
$ws.Range("C91").Value = "Coquimbo"
$ws.Range("D91").Value = 44468
$ws.Range("E91").Value = 5
$ws.Range("F91").Value = 100112010
$ws.Range("G91").Value = "Achicoria"
$ws.Range("H91").Value = "Sin especificar"
$ws.Range("I91").Value = "Primera"
$ws.Range("J91").Value = 140
$ws.Range("K91").Value = 6000
$ws.Range("L91").Value = 6500
$ws.Range("M91").Value = 6286
$ws.Range("N91").Value = "`$/caja 16 unidades"
$ws.Range("O91").Value = "Provincia de Quillota"
$ws.Range("P91").Value = 393
$ws.Range("Q91").Value = 16
$ws.Range("R91").Value = "Hortaliza"
